$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 105: previously a near-empty stub row (only C105 = "taiteien opiskelija"),
# now filled in as a full record for "Fajar Setiawan". ---
$ws.Range("B105").Value = "Fajar Setiawan"
$ws.Range("C105").Value = "musiikin opiskelija"
$ws.Range("D105").Value = "Indonesia"
$ws.Range("E105").Value = 2110
$ws.Range("F105").Formula = "=(2133-E105)"
$ws.Range("G105").Value = 0
$ws.Range("G105").Font.Color = 0
$ws.Range("L105").Value = 1
$ws.Range("N105").Value = "musiikki"
$ws.Range("O105").Value = "löytää paikkansa maailmassa"

# --- Row 107 (Taman Sari): job title changed. ---
$ws.Range("C107").Value = "huoltoteknikko / Pikiran"

# --- New row 108: Dravid Klumm. ---
$ws.Range("A108").Value = 107
$ws.Range("B108").Value = "Dravid Klumm"
$ws.Range("C108").Value = "labrapäällikkö / Pikiran"
$ws.Range("D108").Value = "Intia"
$ws.Range("E108").Value = 2100
$ws.Range("F108").Formula = "=(2133-E108)"
$ws.Range("G108").Value = 1
$ws.Range("G108").Font.Color = 0
$ws.Range("N108").Value = "laboratoriotekniikka"
$ws.Range("O108").Value = "perhe"

# --- Cosmetic: mirror final cursor position left behind in the saved file. ---
$ws.Range("G101").Select() | Out-Null
